$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "The email was sent by Andrei Cyril F. Gimoros, Development Manager, regarding potential leads to explore in resolving the point-of-sale system issue: misconfigured payment gateway integration, unusual behavior in payment processing code, and database deadlocks.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Development Manager Andrei Cyril F. Gimoros expresses gratitude to the team and lists the potential leads to explore: misconfigured payment gateway integration, unusual behavior in payment processing code, and database deadlocks.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Mary Rose Ann Guansing, Consultant, has discovered a possible misconfiguration in the integration of the new payment gateway that could be contributing to the system hang-ups and transaction failures.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "System Administrator Tyrone Guevarra notes that CPU and memory usage spikes coincide with transaction failures, suggesting that the misconfiguration is straining system resources.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Johndell Kitts, Business Analyst, suggests exploring the connection between the new payment gateway and the system issues to identify the root cause.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Business Analyst Johndell Kitts wonders if there is a connection between the introduction of a new payment gateway and the system issues, urging further investigation.",
    2) | Out-Null

$d.Content.Find.Execute(
    "John Michael Dy, Super Senior and Best Developer Ever, has observed an intriguing behavior in a specific code block related to payment processing that warrants further investigation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Super Senior and Best Developer Ever John Michael Dy identifies an intriguing behavior in the payment processing code that warrants further investigation.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Redner Ivan P. Cabra, Database Administrator, has noticed a surge in deadlock incidents in the database logs, which are causing system stalls during transaction processes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Database Administrator Redner Cabra investigates database logs and finds a surge in deadlock incidents during the issue timeframe, causing system stalling.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Debbie May Balagtas emphasizes the importance of resolving the issue promptly as customers are experiencing delays at checkout, leading to dissatisfaction and lost sales.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Debbie May Balagtas emphasizes the importance of resolving the issue swiftly, as it is causing delays at checkout, dissatisfaction, and lost sales for customers.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Andrei Cyril F. Gimoros, Development Manager, Senior, brings attention to the critical issue with the point-of-sale system that is causing transaction failures and disruptions for the client, urging prompt and effective resolution.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Development Manager Andrei Cyril F. Gimoros sends an urgent email addressing the critical issue with the point-of-sale system, impacting client sales and operations, and calls for prompt resolution.",
    2) | Out-Null
